$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "093d69fd815225919ef1a4911c3b54d1f134cc07"
$ws.Range("B1").Value = "EZANA001"
$ws.Range("C1").Value = "Lecturer 001"
$ws.Range("D1").Value = 3456787654
$ws.Range("E1").Value = 2548907654
$ws.Range("F1").Value = "lec001@ezana.org"
$ws.Range("G1").Value = "127001, localhost"
$ws.Range("H1").Value = "79af36419cd79898533a2d9a3028c924fe0d720d"
$ws.Range("I1").Value = "19 Oct 2020"
$ws.Range("J1").Value = "53c904468e7edec9a7f2501d8a8c8d5140c434cb"
